$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.435199999999994
$ws.Range("E3").Value = 13.19699999999999
$ws.Range("B4").Value = 4.649100000000003
$ws.Range("C6").Value = -11.0423
$ws.Range("B7").Value = 6.698500000000001
$ws.Range("C7").Value = -11.7412
$ws.Range("B8").Value = 5.348299999999995
$ws.Range("C8").Value = -10.30469999999999
$ws.Range("E8").Value = 13.08629999999999
$ws.Range("E9").Value = 13.26280000000002
$ws.Range("A11").Value = -21.8095
$ws.Range("D11").Value = -9.00639999999999
$ws.Range("A12").Value = -22.73320000000001
$ws.Range("B12").Value = 6.575600000000001
$ws.Range("B14").Value = 9.153699999999999
$ws.Range("D14").Value = -8.803699999999994
$ws.Range("A15").Value = -21.53320000000003
$ws.Range("C19").Value = -12.12849999999999
$ws.Range("D19").Value = -8.471499999999994
$ws.Range("E19").Value = 13.1354
$ws.Range("C21").Value = -13.31539999999998
$ws.Range("D21").Value = -9.056499999999989
$ws.Range("B22").Value = 5.460400000000002
$ws.Range("C24").Value = -12.0197
$ws.Range("C25").Value = -10.6484

$wb.Save()
